# feat(customer-order): update file template
#
# - Insert a new column before the old "Sub Budget" column, splitting it
#   into "Sub Budget Code" (existing column) + "Sub Budget Name" (new column).
#   This shifts the old "Value"/"Notes" columns one slot to the right.
# - Widen the two budget columns to fit the longer headers.
# - Zoom the sheet view in a bit and move the active selection.
# - Touch the bottom-right corner of the sheet so the saved dimension /
#   row metadata covers the full used range, matching the source template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C; old C (Value) -> D, old D (Notes) -> E.
$ws.Columns.Item(3).EntireColumn.Insert()

# Rename the old "Sub Budget" header and label the newly inserted column.
$ws.Range("B1").Value = "Sub Budget Code"
$ws.Range("C1").Value = "Sub Budget Name"

# Widen both budget columns (character-width units; ~27.49 display width).
$ws.Columns.Item(2).ColumnWidth = 26.65
$ws.Columns.Item(3).ColumnWidth = 26.65

# Bump the zoom level and move the active cell selection to A3.
$ws.Application.ActiveWindow.Zoom = 114
$ws.Range("A3").Select()

# Give the last two rows of the sheet an explicit (small) custom height,
# and touch the bottom-right cell so the sheet's used range / dimension
# extends down to row 1048576 and out to column E.
$ws.Rows.Item(1048575).RowHeight = 12.8
$ws.Rows.Item(1048576).RowHeight = 12.8
$ws.Cells.Item(1048576, 5).NumberFormat = "General"
